# Minijobbeitraege.xlsx - insert the "Zahlt Arbeitnehmer Rentenpauschale?"
# question as a new row right after "kurzfristig beschaeftigt?" and rename
# the two "...beitrag..." rows to "...pauschale...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this pushes the former rows 3-10 down to 4-11
$ws.Rows(3).Insert()

# New row 3: the AN-Rentenpauschale question, defaulting to "ja"
$ws.Cells.Item(3, 1).Value2 = "Zahlt Arbeitnehmer Rentenpauschale?"
$ws.Cells.Item(3, 2).Value2 = "ja"

# Rows 4-6 (formerly 3-5) keep their values but get renamed labels
$ws.Cells.Item(4, 1).Value2 = "Arbeitgeberpauschale Krankenversicherung in Prozent"
$ws.Cells.Item(5, 1).Value2 = "Arbeitgeberpauschale Rentenversicherung in Prozent"
$ws.Cells.Item(6, 1).Value2 = "Arbeitnehmerpauschale Rentenversicherung in Prozent"

# Selection moves to A8, matching the saved cursor position in the edited file
$ws.Range("A8").Select()
